$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -17.60811762772571
$ws.Range("C2").Value = -0.5342166620723915
$ws.Range("D2").Value = -17.60811762772571
$ws.Range("E2").Value = -17.60811762772571
$ws.Range("F2").Value = -17.60811762772571
$ws.Range("G2").Value = -17.60811762772571
$ws.Range("H2").Value = -17.60811762772571
$ws.Range("I2").Value = -17.60811762772571
$ws.Range("J2").Value = -17.60811762772571
$ws.Range("K2").Value = -17.60811762772571

$ws.Range("B3").Value = -17.60811762772571
$ws.Range("C3").Value = -17.60811762772571
$ws.Range("D3").Value = -17.60811762772571
$ws.Range("E3").Value = -17.60811762772571
$ws.Range("F3").Value = -17.60811762772571
$ws.Range("G3").Value = -17.60811762772571
$ws.Range("H3").Value = -17.60811762772571
$ws.Range("I3").Value = -0.4402977002327333
$ws.Range("J3").Value = -17.60811762772571
$ws.Range("K3").Value = -17.60811762772571

$ws.Range("B4").Value = -17.60811762772571
$ws.Range("C4").Value = -0.4258687772048382
$ws.Range("D4").Value = 0.1039861423059535
$ws.Range("E4").Value = -17.60811762772571
$ws.Range("F4").Value = 3.968608893565476
$ws.Range("G4").Value = -17.60811762772571
$ws.Range("H4").Value = 0.881389339085379
$ws.Range("I4").Value = -17.60811762772571
$ws.Range("J4").Value = -17.60811762772571
$ws.Range("K4").Value = -17.60811762772571

$ws.Range("B5").Value = -17.60811762772571
$ws.Range("C5").Value = 0.1754242016191947
$ws.Range("D5").Value = -17.60811762772571
$ws.Range("E5").Value = -17.60811762772571
$ws.Range("F5").Value = -17.60811762772571
$ws.Range("G5").Value = 3.479123116244602
$ws.Range("H5").Value = -17.60811762772571
$ws.Range("I5").Value = -17.60811762772571
$ws.Range("J5").Value = -17.60811762772571
$ws.Range("K5").Value = -17.60811762772571

$ws.Range("B6").Value = -17.60811762772571
$ws.Range("C6").Value = -17.60811762772571
$ws.Range("D6").Value = -17.60811762772571
$ws.Range("E6").Value = -17.60811762772571
$ws.Range("F6").Value = -17.60811762772571
$ws.Range("G6").Value = -17.60811762772571
$ws.Range("H6").Value = -17.60811762772571
$ws.Range("I6").Value = -17.60811762772571
$ws.Range("J6").Value = -17.60811762772571
$ws.Range("K6").Value = -17.60811762772571

$ws.Range("B7").Value = 3.093206603804796
$ws.Range("C7").Value = -17.60811762772571
$ws.Range("D7").Value = -17.60811762772571
$ws.Range("E7").Value = -17.60811762772571
$ws.Range("F7").Value = -17.60811762772571
$ws.Range("G7").Value = -17.60811762772571
$ws.Range("H7").Value = -17.60811762772571
$ws.Range("I7").Value = -17.60811762772571
$ws.Range("J7").Value = -17.60811762772571
$ws.Range("K7").Value = -17.60811762772571

$ws.Range("B8").Value = -17.60811762772571
$ws.Range("C8").Value = -17.60811762772571
$ws.Range("D8").Value = -17.60811762772571
$ws.Range("E8").Value = 2.048329931116285
$ws.Range("F8").Value = -17.60811762772571
$ws.Range("G8").Value = -17.60811762772571
$ws.Range("H8").Value = -17.60811762772571
$ws.Range("I8").Value = -17.60811762772571
$ws.Range("J8").Value = -17.60811762772571
$ws.Range("K8").Value = -17.60811762772571

$ws.Range("B9").Value = 3.519290484345814
$ws.Range("C9").Value = -17.60811762772571
$ws.Range("D9").Value = -17.60811762772571
$ws.Range("E9").Value = -17.60811762772571
$ws.Range("F9").Value = -17.60811762772571
$ws.Range("G9").Value = -17.60811762772571
$ws.Range("H9").Value = -17.60811762772571
$ws.Range("I9").Value = -17.60811762772571
$ws.Range("J9").Value = -17.60811762772571
$ws.Range("K9").Value = -17.60811762772571

$ws.Range("B10").Value = -17.60811762772571
$ws.Range("C10").Value = -17.60811762772571
$ws.Range("D10").Value = -17.60811762772571
$ws.Range("E10").Value = -17.60811762772571
$ws.Range("F10").Value = -17.60811762772571
$ws.Range("G10").Value = -17.60811762772571
$ws.Range("H10").Value = -17.60811762772571
$ws.Range("I10").Value = 0.8604626666090643
$ws.Range("J10").Value = -17.60811762772571
$ws.Range("K10").Value = 2.068715814847168

$ws.Range("B11").Value = -17.60811762772571
$ws.Range("C11").Value = -17.60811762772571
$ws.Range("D11").Value = -17.60811762772571
$ws.Range("E11").Value = 1.972930286597393
$ws.Range("F11").Value = -17.60811762772571
$ws.Range("G11").Value = 1.389053581866031
$ws.Range("H11").Value = -17.60811762772571
$ws.Range("I11").Value = -17.60811762772571
$ws.Range("J11").Value = -17.60811762772571
$ws.Range("K11").Value = 1.365033176650472

$ws.Range("B12").Value = -17.60811762772571
$ws.Range("C12").Value = -17.60811762772571
$ws.Range("D12").Value = -17.60811762772571
$ws.Range("E12").Value = -17.60811762772571
$ws.Range("F12").Value = -17.60811762772571
$ws.Range("G12").Value = -17.60811762772571
$ws.Range("H12").Value = -17.60811762772571
$ws.Range("I12").Value = -17.60811762772571
$ws.Range("J12").Value = -17.60811762772571
$ws.Range("K12").Value = -17.60811762772571

$ws.Range("B13").Value = -17.60811762772571
$ws.Range("C13").Value = -17.60811762772571
$ws.Range("D13").Value = -17.60811762772571
$ws.Range("E13").Value = 1.784730124821259
$ws.Range("F13").Value = -17.60811762772571
$ws.Range("G13").Value = -17.60811762772571
$ws.Range("H13").Value = -17.60811762772571
$ws.Range("I13").Value = -17.60811762772571
$ws.Range("J13").Value = -17.60811762772571
$ws.Range("K13").Value = 1.689615902258428

$ws.Range("B14").Value = -17.60811762772571
$ws.Range("C14").Value = -17.60811762772571
$ws.Range("D14").Value = 1.770915881696862
$ws.Range("E14").Value = -17.60811762772571
$ws.Range("F14").Value = -17.60811762772571
$ws.Range("G14").Value = -17.60811762772571
$ws.Range("H14").Value = -17.60811762772571
$ws.Range("I14").Value = -17.60811762772571
$ws.Range("J14").Value = -17.60811762772571
$ws.Range("K14").Value = 1.328178113669579

$ws.Range("B15").Value = -17.60811762772571
$ws.Range("C15").Value = -17.60811762772571
$ws.Range("D15").Value = -0.4850773685743109
$ws.Range("E15").Value = -17.60811762772571
$ws.Range("F15").Value = -17.60811762772571
$ws.Range("G15").Value = -17.60811762772571
$ws.Range("H15").Value = -17.60811762772571
$ws.Range("I15").Value = -17.60811762772571
$ws.Range("J15").Value = -17.60811762772571
$ws.Range("K15").Value = -17.60811762772571

$ws.Range("B16").Value = -17.60811762772571
$ws.Range("C16").Value = -17.60811762772571
$ws.Range("D16").Value = -17.60811762772571
$ws.Range("E16").Value = -17.60811762772571
$ws.Range("F16").Value = -17.60811762772571
$ws.Range("G16").Value = -17.60811762772571
$ws.Range("H16").Value = -17.60811762772571
$ws.Range("I16").Value = -17.60811762772571
$ws.Range("J16").Value = 4.321921234832267
$ws.Range("K16").Value = -17.60811762772571

$ws.Range("B17").Value = -17.60811762772571
$ws.Range("C17").Value = 0.06019261211739748
$ws.Range("D17").Value = -0.4536851589348114
$ws.Range("E17").Value = -17.60811762772571
$ws.Range("F17").Value = -17.60811762772571
$ws.Range("G17").Value = -17.60811762772571
$ws.Range("H17").Value = 3.053748305069607
$ws.Range("I17").Value = -0.6756144947809635
$ws.Range("J17").Value = -17.60811762772571
$ws.Range("K17").Value = -17.60811762772571

$ws.Range("B18").Value = -17.60811762772571
$ws.Range("C18").Value = -17.60811762772571
$ws.Range("D18").Value = -17.60811762772571
$ws.Range("E18").Value = -17.60811762772571
$ws.Range("F18").Value = -17.60811762772571
$ws.Range("G18").Value = -17.60811762772571
$ws.Range("H18").Value = 2.072944019468157
$ws.Range("I18").Value = -0.4919666073578666
$ws.Range("J18").Value = -17.60811762772571
$ws.Range("K18").Value = -17.60811762772571

$ws.Range("B19").Value = -17.60811762772571
$ws.Range("C19").Value = -17.60811762772571
$ws.Range("D19").Value = 2.939696644222277
$ws.Range("E19").Value = -17.60811762772571
$ws.Range("F19").Value = -17.60811762772571
$ws.Range("G19").Value = -17.60811762772571
$ws.Range("H19").Value = 0.967218489466767
$ws.Range("I19").Value = 2.61175161763022
$ws.Range("J19").Value = -17.60811762772571
$ws.Range("K19").Value = -17.60811762772571

$ws.Range("B20").Value = -17.60811762772571
$ws.Range("C20").Value = 3.074246935091408
$ws.Range("D20").Value = 2.677023788948343
$ws.Range("E20").Value = -17.60811762772571
$ws.Range("F20").Value = 2.119119822675309
$ws.Range("G20").Value = -17.60811762772571
$ws.Range("H20").Value = 0.879349534021393
$ws.Range("I20").Value = 3.321603831352709
$ws.Range("J20").Value = -17.60811762772571
$ws.Range("K20").Value = 2.905429751149777

$ws.Range("B21").Value = -17.60811762772571
$ws.Range("C21").Value = 2.994690863163727
$ws.Range("D21").Value = -17.60811762772571
$ws.Range("E21").Value = 3.086183290908284
$ws.Range("F21").Value = -17.60811762772571
$ws.Range("G21").Value = 2.639157435167036
$ws.Range("H21").Value = 0.8891299978205868
$ws.Range("I21").Value = -17.60811762772571
$ws.Range("J21").Value = -17.60811762772571
$ws.Range("K21").Value = -17.60811762772571

